$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.370.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.97%  "

$ws.Range("D3").Value = "'1.827.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'314.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.66%  "

$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").Value = "'0.4477"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.28%  "

$ws.Range("D8").Value = "'0.3782"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.34%  "

$ws.Range("D9").Value = "'0.07481"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.09%  "

$ws.Range("D10").Value = "'0.8874"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.12%  "

$ws.Range("D11").Value = "'21.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "

$ws.Range("D12").Value = "'1.829.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.34%  "

$ws.Range("D13").Value = "'6.764"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.04%  "

$ws.Range("D14").Value = "'5.444"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.85%  "

$ws.Range("D15").Value = "'93.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.95%  "

$ws.Range("D16").Value = "'0.07122"
$ws.Range("D16").Style = "Normal"

$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").Value = "'0.000008798"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.39%  "

$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").Value = "'15.17"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").Value = "'27.381.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.21%  "

$ws.Range("D22").Value = "'5.409"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.08%  "

$ws.Range("E23").Value = "  -0.46%  "

$ws.Range("D24").Value = "'2.057.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.78%  "

$ws.Range("D25").Value = "'1.968"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.70%  "

$ws.Range("D26").Value = "'151.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").Value = "'2.314"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.87%  "

$ws.Range("E28").Value = "  +0.92%  "

$ws.Range("D29").Value = "'5.404"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.53%  "

$ws.Range("D30").Value = "'117.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.30%  "

$ws.Range("D31").Value = "'0.08888"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.26%  "

$ws.Range("D32").Value = "'0.7918"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.85%  "

$ws.Range("D33").Value = "'1.209"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.13%  "

$ws.Range("E34").Value = "  +2.89%  "

$ws.Range("D35").Value = "'2.925"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.40%  "

$ws.Range("D36").Value = "'0.9997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  +0.52%  "

$ws.Range("D38").Value = "'0.01992"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.42%  "

$ws.Range("D39").Value = "'0.05312"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("D40").Value = "'7.340"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.32%  "

$ws.Range("D41").Value = "'0.5361"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.14%  "

$ws.Range("D42").Value = "'2.868"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.75%  "

$ws.Range("D43").Value = "'0.1723"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.65%  "

$ws.Range("D44").Value = "'2.332"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +17.66%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'8.694"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.86%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5137"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.57%  "

$ws.Range("D47").Value = "'10.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.82%  "

$ws.Range("D48").Value = "'1.701"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.44%  "

$ws.Range("D49").Value = "'105.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.60%  "

$ws.Range("D50").Value = "'0.9998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("D51").Value = "'0.06412"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.46%  "
